# ----------------------------------------------------------------------------
# Update Liga_brasil_B_2025 sheet: add two possession columns (V, W), fix
# first/second-half goal tallies on several existing fixtures, and append the
# newest round of matches (rows 201-211).
# ----------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) New header columns V1/W1 -- reuse U1s (bold/bordered/centered) style.
$ws.Range("U1").Copy()
$ws.Range("V1:W1").PasteSpecial(-4122)
$ws.Range("V1").Value = "Posesión Local ().2"
$ws.Range("W1").Value = "Posesión Visita ().2"

# 2) Corrected first/second-half goal splits on existing fixtures (M/N/O/P).
$ws.Range("M159").Value = 2
$ws.Range("O159").Value = 2
$ws.Range("M161").Value = 2
$ws.Range("N161").Value = 1
$ws.Range("O161").Value = 1
$ws.Range("P161").Value = 0
$ws.Range("M162").Value = 1
$ws.Range("N162").Value = 1
$ws.Range("O162").Value = 0
$ws.Range("P162").Value = 0
$ws.Range("M163").Value = 1
$ws.Range("N163").Value = 1
$ws.Range("O163").Value = 0
$ws.Range("P163").Value = 0
$ws.Range("M166").Value = 1
$ws.Range("O166").Value = 2
$ws.Range("M168").Value = 1
$ws.Range("N168").Value = 3
$ws.Range("O168").Value = 1
$ws.Range("P168").Value = 2
$ws.Range("N169").Value = 2
$ws.Range("P169").Value = 0
$ws.Range("N171").Value = 1
$ws.Range("P171").Value = 0
$ws.Range("M174").Value = 1
$ws.Range("N174").Value = 1
$ws.Range("O174").Value = 0
$ws.Range("P174").Value = 0
$ws.Range("M175").Value = 1
$ws.Range("O175").Value = 1
$ws.Range("M176").Value = 3
$ws.Range("N176").Value = 1
$ws.Range("O176").Value = 1
$ws.Range("P176").Value = 1
$ws.Range("M177").Value = 1
$ws.Range("O177").Value = 0
$ws.Range("M178").Value = 1
$ws.Range("N178").Value = 1
$ws.Range("O178").Value = 2
$ws.Range("P178").Value = 0
$ws.Range("N179").Value = 1
$ws.Range("P179").Value = 0
$ws.Range("M180").Value = 2
$ws.Range("N180").Value = 1
$ws.Range("O180").Value = 0
$ws.Range("P180").Value = 0
$ws.Range("M182").Value = 1
$ws.Range("O182").Value = 3
$ws.Range("M183").Value = 1
$ws.Range("O183").Value = 0
$ws.Range("M184").Value = 1
$ws.Range("O184").Value = 1
$ws.Range("N185").Value = 1
$ws.Range("P185").Value = 1
$ws.Range("M187").Value = 4
$ws.Range("O187").Value = 1
$ws.Range("M188").Value = 1
$ws.Range("N188").Value = 1
$ws.Range("O188").Value = 0
$ws.Range("P188").Value = 0
$ws.Range("M189").Value = 1
$ws.Range("O189").Value = 0
$ws.Range("N190").Value = 1
$ws.Range("P190").Value = 0
$ws.Range("M191").Value = 1
$ws.Range("N191").Value = 2
$ws.Range("O191").Value = 0
$ws.Range("P191").Value = 0
$ws.Range("N192").Value = 1
$ws.Range("P192").Value = 1
$ws.Range("M194").Value = 1
$ws.Range("N194").Value = 1
$ws.Range("O194").Value = 1
$ws.Range("P194").Value = 0
$ws.Range("M195").Value = 1
$ws.Range("N195").Value = 2
$ws.Range("O195").Value = 1
$ws.Range("P195").Value = 0
$ws.Range("M196").Value = 3
$ws.Range("O196").Value = 0
$ws.Range("M197").Value = 1
$ws.Range("N197").Value = 1
$ws.Range("O197").Value = 0
$ws.Range("P197").Value = 0
$ws.Range("M199").Value = 1
$ws.Range("N199").Value = 1
$ws.Range("O199").Value = 0
$ws.Range("P199").Value = 0
$ws.Range("M200").Value = 1
$ws.Range("O200").Value = 1

# 3) Append the 11 newest fixtures (rows 201-211).
$newRows = @(
    ,@(201, "2025-08-08", "Ferroviária", "Amazonas", 2, 1, 1353468, 9, 1, 2, 3, 0, 0, 1, 0, 1, 1, 60, 40, "L")
    ,@(202, "2025-08-09", "Coritiba", "Chapecoense-sc", 0, 0, 1353460, 4, 3, 3, 3, 0, 0, 0, 0, 0, 0, 54, 46, "E")
    ,@(203, "2025-08-09", "America Mineiro", "remo", 0, 1, 1353467, 9, 3, 1, 1, 0, 0, 0, 1, 0, 0, 61, 39, "V")
    ,@(204, "2025-08-09", "Goias", "Operario-PR", 2, 1, 1353461, 4, 3, 1, 4, 0, 0, 1, 1, 1, 0, 43, 57, "L")
    ,@(205, "2025-08-09", "Volta Redonda", "Novorizontino", 0, 0, 1353469, 11, 4, 3, 0, 0, 0, 0, 0, 0, 0, 53, 47, "E")
    ,@(206, "2025-08-10", "Avai", "Cuiaba", 2, 0, 1353464, 1, 10, 1, 2, 0, 0, 2, 0, 0, 0, 36, 64, "L")
    ,@(207, "2025-08-10", "Atletico Goianiense", "Botafogo SP", 2, 0, 1353462, 5, 5, 2, 2, 0, 0, 0, 0, 2, 0, 61, 39, "L")
    ,@(208, "2025-08-11", "Criciuma", "Atletico Paranaense", 4, 2, 1353465, 4, 5, 3, 3, 0, 0, 1, 2, 3, 0, 51, 49, "L")
    ,@(209, "2025-08-12", "Paysandu", "Vila Nova", 0, 1, 1353463, 9, 5, 1, 2, 0, 0, 0, 0, 0, 1, 59, 41, "V")
    ,@(210, "2025-08-12", "CRB", "Athletic Club", 1, 0, 1353466, 8, 4, 1, 2, 0, 1, 0, 0, 1, 0, 48, 52, "L")
    ,@(211, "2025-08-15", "Novorizontino", "Coritiba", 1, 2, 1353478, 4, 6, 5, 4, 0, 0, 0, 1, 1, 1, 58, 42, "V")
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
    $ws.Cells.Item($rowNum, 7).Value = $r[7]
    $ws.Cells.Item($rowNum, 8).Value = $r[8]
    $ws.Cells.Item($rowNum, 9).Value = $r[9]
    $ws.Cells.Item($rowNum, 10).Value = $r[10]
    $ws.Cells.Item($rowNum, 11).Value = $r[11]
    $ws.Cells.Item($rowNum, 12).Value = $r[12]
    $ws.Cells.Item($rowNum, 13).Value = $r[13]
    $ws.Cells.Item($rowNum, 14).Value = $r[14]
    $ws.Cells.Item($rowNum, 15).Value = $r[15]
    $ws.Cells.Item($rowNum, 16).Value = $r[16]
    $ws.Cells.Item($rowNum, 17).Value = $r[17]
    $ws.Cells.Item($rowNum, 18).Value = $r[18]
    $ws.Cells.Item($rowNum, 19).Value = $r[19]
}

